$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate existing entries (re-creates shared strings table in new order)
$ws.Range("A1").Value = "15.03.2022 - 14.04.2022"
$ws.Range("B2").Value = "food"
$ws.Range("B3").Value = "whater"
$ws.Range("B4").Value = "medicine"

# Add new expense row for guests (18.03.2022)
$ws.Range("B5").Value = "guests"
$ws.Range("C5").Value = 80

# Update the active selection to match the authored workbook state
$ws.Range("D10").Select() | Out-Null
